$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell carrying the default (unstyled) cell format used throughout column D,
# so that forcing text below does not leave a stray style on the edited cells.
$defaultStyle = $ws.Range("D19").Style

# Map of cell -> new value for this update. Column D values are prefixed with a single
# quote to force Excel to keep them as text (otherwise values like '14.77' or '307.60'
# get auto-converted to numbers and lose their original formatting).
$updates = @(
    ,@('D2', '''27.280.68')
    ,@('E2', '  +0.36%  ')
    ,@('D3', '''1.906.98')
    ,@('E3', '  +0.28%  ')
    ,@('D4', '''1.002')
    ,@('E4', '  +0.11%  ')
    ,@('D5', '''307.60')
    ,@('E5', '  -0.10%  ')
    ,@('D6', '''1.001')
    ,@('E6', '  +0.18%  ')
    ,@('D7', '''0.5268')
    ,@('E7', '  +1.06%  ')
    ,@('D8', '''0.3817')
    ,@('E8', '  +1.35%  ')
    ,@('D9', '''0.07297')
    ,@('E9', '  +0.26%  ')
    ,@('D10', '''22.11')
    ,@('E10', '  +4.56%  ')
    ,@('D11', '''0.9019')
    ,@('E11', '  -0.29%  ')
    ,@('D12', '''0.08197')
    ,@('E12', '  -1.26%  ')
    ,@('D13', '''95.78')
    ,@('E13', '  -0.98%  ')
    ,@('D14', '''5.352')
    ,@('E14', '  +1.20%  ')
    ,@('D15', '''1.001')
    ,@('E15', '  +0.08%  ')
    ,@('B16', 'ShibaInu')
    ,@('C16', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib')
    ,@('D16', '''0.000008642')
    ,@('E16', '  -0.28%  ')
    ,@('B17', 'Avalanche')
    ,@('C17', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax')
    ,@('D17', '''14.77')
    ,@('E17', '  +1.47%  ')
    ,@('D18', '''1.333.61')
    ,@('E18', '  -30.37%  ')
    ,@('E19', '  +0.21%  ')
    ,@('D20', '''27.313.06')
    ,@('E20', '  +0.33%  ')
    ,@('D21', '''5.070')
    ,@('E21', '  -0.38%  ')
    ,@('D22', '''10.82')
    ,@('E22', '  +1.66%  ')
    ,@('D23', '''6.518')
    ,@('E23', '  +1.23%  ')
    ,@('D24', '''149.92')
    ,@('E24', '  +2.40%  ')
    ,@('D25', '''2.297')
    ,@('E25', '  -1.04%  ')
    ,@('D26', '''18.24')
    ,@('E26', '  -0.06%  ')
    ,@('D27', '''1.736')
    ,@('E27', '  -0.60%  ')
    ,@('D28', '''116.31')
    ,@('E28', '  +1.12%  ')
    ,@('D29', '''4.835')
    ,@('E29', '  -0.07%  ')
    ,@('D30', '''4.821')
    ,@('E30', '  -1.49%  ')
    ,@('D31', '''0.09260')
    ,@('E31', '  -0.10%  ')
    ,@('D32', '''0.8384')
    ,@('E32', '  +5.02%  ')
    ,@('D33', '''0.05073')
    ,@('E33', '  -0.21%  ')
    ,@('D34', '''1.227')
    ,@('E34', '  -1.46%  ')
    ,@('E35', '  +1.76%  ')
    ,@('E36', '  -1.97%  ')
    ,@('D37', '''2.683')
    ,@('E37', '  +3.39%  ')
    ,@('D38', '''0.5753')
    ,@('E38', '  +0.53%  ')
    ,@('D39', '''0.02004')
    ,@('E39', '  +0.13%  ')
    ,@('D40', '''1.077')
    ,@('E40', '  -0.20%  ')
    ,@('D41', '''9.316')
    ,@('E41', '  +3.26%  ')
    ,@('D42', '''6.531')
    ,@('E42', '  -1.02%  ')
    ,@('D43', '''116.84')
    ,@('E43', '  -0.22%  ')
    ,@('E44', '  +0.20%  ')
    ,@('E45', '  +1.12%  ')
    ,@('E46', '  +0.17%  ')
    ,@('D47', '''10.16')
    ,@('E47', '  +0.85%  ')
    ,@('D48', '''1.636')
    ,@('E48', '  +0.32%  ')
    ,@('D49', '''38.85')
    ,@('E49', '  +2.95%  ')
    ,@('D50', '''0.06178')
    ,@('E50', '  +3.73%  ')
    ,@('D51', '''63.86')
    ,@('E51', '  -0.22%  ')
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# Column D cells that were forced to text above: strip the induced quote-prefix/number
# style so the cell format matches the sheet default again.
$dCells = @('D2', 'D3', 'D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($c in $dCells) {
    $ws.Range($c).Style = $defaultStyle
}
